$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Shape "Rectangle 5" (first TWAMP diagram box) ---------------------
$r5 = $s.Shapes.Item(3)
$tr5 = $r5.TextFrame.TextRange

# Paragraph 4: "  |         Transmit Timestamp                                    |"
# -> split into 3 runs, bolding the "(t1)" label portion.
$p5t = $tr5.Paragraphs(4, 1)
$p5t.Characters(3, 64).Text = "|         Transmit Timestamp (t1)                               "
$p5t.Characters(3, 64).Font.Bold = $true

# Paragraph 9: "  |         Receive Timestamp                                     |"
# -> split into 3 runs, bolding the "(t2)" label portion.
$p5r = $tr5.Paragraphs(9, 1)
$p5r.Characters(3, 64).Text = "|         Receive Timestamp (t2)                                "
$p5r.Characters(3, 64).Font.Bold = $true

# --- Shape "Rectangle 6" (second, STAMP diagram box) --------------------
$r6 = $s.Shapes.Item(6)
$tr6 = $r6.TextFrame.TextRange

# Paragraph 4: "  |         Transmit Timestamp                                    |"
# -> whole run becomes bold and gains the "(t1)" suffix.
$p6t = $tr6.Paragraphs(4, 1)
$p6t.Characters(1, 67).Text = "  |         Transmit Timestamp (t1)                               |"
$p6t.Characters(1, 67).Font.Bold = $true

# Paragraph 9: "  |         Receive Timestamp                                     |"
# -> whole run becomes bold and gains the "(t2)" suffix.
$p6r = $tr6.Paragraphs(9, 1)
$p6r.Characters(1, 67).Text = "  |         Receive Timestamp (t2)                                |"
$p6r.Characters(1, 67).Font.Bold = $true

# --- Handout master footer date field -----------------------------------
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Value = "7/9/20"
